# Recipe.docx - "added line to recipe"
#
# The second paragraph's single run of body text is reworked into four
# runs (the middle phrase "all of" becomes its own run, and a brand new
# closing sentence is appended as a further run) while keeping the
# existing formatting (Century Schoolbook / 24 half-points) identical
# throughout.

$d = $word.ActiveDocument

# 1. Append the new closing sentence right after "screen for." using
#    Find/Replace on the whole document Content range. Doing the
#    replacement this way inserts the new text *inside* the paragraph,
#    ahead of the trailing hidden "_GoBack" bookmark, exactly where the
#    diff places it.
$tail = $d.Content
$tail.Find.Execute(
    "screen for.", $true, $false, $false, $false, $false, $true, 1, $false,
    "screen for. The function works by having the if else statement run and then in the else code the number of the image per the variable has 1 added to it so image1 will go to image2.",
    2
)

# 2. Give the freshly appended sentence its own run: flip a character
#    format on and straight back off. That forces a run boundary at the
#    sentence's edges without changing how the text actually looks.
$newSentence = $d.Content
$newSentence.Find.Execute(" The function works by having the if else statement run and then in the else code the number of the image per the variable has 1 added to it so image1 will go to image2.")
$newSentence.Font.Bold = $true
$newSentence.Font.Bold = $false

# 3. Split "all of" out of the first sentence into its own run the same
#    way (this is where Word's grammar checker would normally have
#    flagged the phrase while the paragraph was being retyped).
$allOf = $d.Content
$allOf.Find.Execute("all of")
$allOf.Font.Bold = $true
$allOf.Font.Bold = $false
